# Actualización automática de catálogo y fotos
#
# Applies the content changes described in the commit diff to the
# "catalogo" worksheet of the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("catalogo")

# --- Swap a few image filename pairs (imagen1/imagen2 columns) ---
# (use Value2 for reads -- Value reads are unreliable in this runtime)

# Chicago Terciopelo Negro (row 2): swap E2/F2
$tmp = $ws.Range("E2").Value2
$ws.Range("E2").Value = $ws.Range("F2").Value2
$ws.Range("F2").Value = $tmp

# Chicago Blanco Cocodrilo (row 5): swap E5/F5
$tmp = $ws.Range("E5").Value2
$ws.Range("E5").Value = $ws.Range("F5").Value2
$ws.Range("F5").Value = $tmp

# Chicago Naranja (row 6): swap E6/F6
$tmp = $ws.Range("E6").Value2
$ws.Range("E6").Value = $ws.Range("F6").Value2
$ws.Range("F6").Value = $tmp

# San Francisco Azul (row 25): swap E25/F25
$tmp = $ws.Range("E25").Value2
$ws.Range("E25").Value = $ws.Range("F25").Value2
$ws.Range("F25").Value = $tmp

# --- Rename a photo file extension ---
# Vancouver Rafia (row 35): imagen1 .png -> .jpg
$ws.Range("E35").Value = "vancouver rafia verde 1.jpg"

# --- Clear a handful of stray formatted-but-empty cells ---
$ws.Range("H3").Clear()
$ws.Range("H5").Clear()
$ws.Range("G7").Clear()
$ws.Range("F24").Clear()

# --- Paris Gorro: split into three colour variants ---
# Row 42 stays as the "red" hat
$ws.Range("A42").Value = "Paris Gorro Rojo"

# Row 43 becomes the "green" hat, and gets the two existing green photos
# (swapped order versus before)
$ws.Range("A43").Value = "Paris Gorro Verde"
$ws.Range("E43").Value = "gorro paris verde.jpeg"
$ws.Range("F43").Value = "gorro paris verde cerca.jpeg"

# Row 44 becomes a brand new "lila" hat variant (no photos yet)
$ws.Range("A44").Value = "Paris Gorro Lila"
